# Update the "取得日時" (retrieved timestamp) column on sheet "案件情報"
# for the existing rows (2-9) to reflect the latest append run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2026-01-07 12:54:36"

for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
